$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.071.14'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -3.64%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.520.07'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -4.54%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '580.95'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.40%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '174.90'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -2.71%  '
$ws.Range('E7').Value = '  +0.48%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.513.15'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -4.54%  '
$ws.Range('E9').Value = '  +0.07%  '
$ws.Range('E10').Value = '  -5.63%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.78'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +7.44%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.602'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -1.80%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '47.35'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -5.21%  '
$ws.Range('E14').Value = '  -3.15%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '675.75'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -1.37%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.087.52'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -4.51%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '8.80'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -2.19%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.516.96'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -4.34%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '69.074.66'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -3.64%  '
$ws.Range('E20').Value = '  -1.31%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.59'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -2.79%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '11.28'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -3.23%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.910'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -3.25%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '16.24'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -8.82%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '98.30'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -5.44%  '
$ws.Range('E26').Value = '  -4.15%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '5.83'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.34%  '
$ws.Range('E28').Value = '  -5.85%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.50'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -6.93%  '
$ws.Range('E31').Value = '  -6.85%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.79'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -4.73%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.22'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -7.89%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '7.41'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +0.86%  '
$ws.Range('E35').Value = '  -5.08%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '579.55'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +1.14%  '
$ws.Range('E37').Value = '  -15.26%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '10.94'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -3.32%  '
$ws.Range('E39').Value = '  -3.39%  '
$ws.Range('E40').Value = '  -3.52%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.00'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +0.21%  '
$ws.Range('E42').Value = '  -3.44%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0442'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -4.47%  '
$ws.Range('E44').Value = '  -6.09%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.439.50'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -8.99%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '33.61'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -5.19%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0₃0709'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -8.84%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.92'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +1.37%  '
$ws.Range('E49').Value = '  -6.56%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '131.99'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -1.93%  '
